$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data values in row 2 -----------------------------------------------
# Row 2 cells carry cell style s="1", whose number format is the built-in Text format
# (numFmtId 49, "@"). Writing a number straight into a Text-formatted cell makes Excel
# store it as a text string, which is not what the workbook needs here (the refreshed
# values must stay real numbers while the cells keep their existing Text style). So for
# every touched cell we momentarily flip its number format to a numeric one, write the
# new number, then flip the format back to Text - the end result is a plain numeric
# value under the same s="1" style, matching how the source workbook was edited.
$updatedValues = [ordered]@{
    "B2" = 2.5
    "E2" = 2
    "F2" = 182
    "M2" = 1.44
    "N2" = 1.89
    "O2" = 1.8
    "P2" = 1.63
    "Q2" = 1.6
    "AS2" = 0.29
    "AT2" = 0.17
    "AU2" = 14.44
    "AV2" = 10.84
    "AX2" = 106.47
    "EK2" = 0.07
    "EL2" = 13.66
    "EM2" = 7.27
    "EN2" = 49.11
    "EO2" = 148.96
    "EP2" = 0.33
    "EQ2" = 0.16
    "ER2" = 0.16
    "ES2" = 20.08
    "ET2" = 11.21
    "EV2" = 145.47
    "EW2" = -44.9
    "EX2" = -56.25
    "EY2" = -31.96
    "EZ2" = -35.13
    "FA2" = 2.14
    "FB2" = 2.4
    "FH2" = 26
    "FI2" = 4889
    "FJ2" = 2603.99
    "FK2" = 17582.12
    "FL2" = 53326.91
    "FN2" = 2005.99
    "FO2" = 4227
    "FP2" = 17582.12
    "FQ2" = 182
}

foreach ($addr in $updatedValues.Keys) {
    $ws.Range($addr).NumberFormat = "0"
}
foreach ($addr in $updatedValues.Keys) {
    $ws.Range($addr).Value = $updatedValues[$addr]
}
foreach ($addr in $updatedValues.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Update the cursor / active-cell selection on the sheet -----------------------
$ws.Range("P13").Select()

